# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 160 (pushing all existing
# rows 160..241 down to 161..242, so the former row 241 becomes row 242).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 160; Excel shifts rows
# 160..241 down to 161..242 and extends the used range accordingly.
$ws.Rows.Item(160).Insert()

# Populate the newly inserted row 160 with the new weekly record.
$ws.Cells.Item(160, 1).Value = 4
$ws.Cells.Item(160, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(160, 3).Value = "Los Lagos"
$ws.Cells.Item(160, 4).Value = 44719
$ws.Cells.Item(160, 5).Value = 10
$ws.Cells.Item(160, 6).Value = "Fruta"
$ws.Cells.Item(160, 7).Value = 100109
$ws.Cells.Item(160, 8).Value = "Uva"
$ws.Cells.Item(160, 9).Value = 100109001
$ws.Cells.Item(160, 10).Value = "Uva"
$ws.Cells.Item(160, 11).Value = "Red Globe"
$ws.Cells.Item(160, 12).Value = "Primera"
$ws.Cells.Item(160, 13).Value = 200
$ws.Cells.Item(160, 14).Value = 11000
$ws.Cells.Item(160, 15).Value = 12000
$ws.Cells.Item(160, 16).Value = 11500
$ws.Cells.Item(160, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(160, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(160, 19).Value = 575
$ws.Cells.Item(160, 20).Value = 20
